# Removed old userManagement tests from del2
#
# - Remove the old "t1" test-fixture row (row 2): clear its contents so the
#   row disappears from the sheet (rows 3-10 keep their original row
#   numbers, dimension becomes A3:...).
# - Row 3 ("t2") gains a Department value in column H ("Biology").
# - A new test-fixture row 11 is appended with a new user
#   (newuser@example.com / password / Student / false / "" / "" / 0).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old row 2 fixture entirely.
$ws.Rows("2:2").ClearContents()

# Row 3 gets a new "Biology" department tag in column H.
$ws.Cells.Item(3, 8).Value = "Biology"

# Append the new user fixture as row 11.
$ws.Cells.Item(11, 2).Value = 5
$ws.Cells.Item(11, 3).Value = "newuser@example.com"
$ws.Cells.Item(11, 4).Value = "password"
$ws.Cells.Item(11, 5).Value = "Student"
$ws.Cells.Item(11, 6).Value = $false
$ws.Cells.Item(11, 7).Value = ""
$ws.Cells.Item(11, 8).Value = ""
$ws.Cells.Item(11, 9).Value = 0
